$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.614.95'
$ws.Range("E2").Value = '  +2.87%  '

$ws.Range("D3").Value = '2.322.89'
$ws.Range("E3").Value = '  +2.04%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = "'516.13"
$ws.Range("E5").Value = '  +1.86%  '

$ws.Range("D6").Value = "'135.52"
$ws.Range("E6").Value = '  +5.43%  '

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = "'0.537"
$ws.Range("E8").Value = '  +1.46%  '

$ws.Range("D9").Value = '2.343.57'
$ws.Range("E9").Value = '  +2.40%  '

$ws.Range("E10").Value = '  +3.59%  '

$ws.Range("E11").Value = '  -1.24%  '

$ws.Range("D12").Value = "'5.36"
$ws.Range("E12").Value = '  +5.57%  '

$ws.Range("E13").Value = '  +0.55%  '

$ws.Range("D14").Value = "'23.93"
$ws.Range("E14").Value = '  +1.09%  '

$ws.Range("D15").Value = '2.739.71'
$ws.Range("E15").Value = '  +2.12%  '

$ws.Range("D16").Value = '56.669.34'
$ws.Range("E16").Value = '  +2.90%  '

$ws.Range("E17").Value = '  +2.50%  '

$ws.Range("D18").Value = '2.320.22'
$ws.Range("E18").Value = '  +2.29%  '

$ws.Range("E19").Value = '  +1.56%  '

$ws.Range("D20").Value = "'326.19"
$ws.Range("E20").Value = '  +3.70%  '

$ws.Range("D21").Value = "'4.22"
$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("D22").Value = "'6.60"
$ws.Range("E22").Value = '  +0.72%  '

$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").Value = "'60.80"
$ws.Range("E24").Value = '  +1.56%  '

$ws.Range("E25").Value = '  +6.26%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("E27").Value = '  +5.65%  '

$ws.Range("E28").Value = '  +11.56%  '

$ws.Range("E29").Value = '  +5.17%  '

$ws.Range("D30").Value = "'169.31"
$ws.Range("E30").Value = '  -1.02%  '

$ws.Range("D31").Value = "'1.69"
$ws.Range("E31").Value = '  +3.11%  '

$ws.Range("E32").Value = '  +0.87%  '

$ws.Range("D33").Value = "'18.48"
$ws.Range("E33").Value = '  +2.87%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("E36").Value = '  +2.15%  '

$ws.Range("D37").Value = "'0.917"
$ws.Range("E37").Value = '  +1.41%  '

$ws.Range("D38").Value = "'4.01"
$ws.Range("E38").Value = '  +3.07%  '

$ws.Range("E39").Value = '  +7.04%  '

$ws.Range("D40").Value = "'38.34"
$ws.Range("E40").Value = '  +4.19%  '

$ws.Range("E41").Value = '  +1.78%  '

$ws.Range("D42").Value = "'142.17"
$ws.Range("E42").Value = '  +3.95%  '

$ws.Range("D43").Value = "'3.61"
$ws.Range("E43").Value = '  +3.83%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = "'278.27"
$ws.Range("E44").Value = '  +7.53%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = "'5.18"
$ws.Range("E45").Value = '  +5.97%  '

$ws.Range("E46").Value = '  +1.78%  '

$ws.Range("E47").Value = '  +0.12%  '

$ws.Range("E48").Value = '  +2.58%  '

$ws.Range("E49").Value = '  +2.59%  '

$ws.Range("E50").Value = '  +8.87%  '

$ws.Range("E51").Value = '  +1.63%  '
